# ADD CORRECCION DE PU, REPROCESOS Y WF
# This script edits the "08" worksheet (2nd sheet) of the workbook:
#  - Splits the QA verification tasks (rows 2-7) into Pago Unico / Reprocesos / Workflow variants
#  - Adds two new tasks (rows 8-9) for "Verificacion del Control de Calidad"
#  - Adjusts merged cells, row heights, column F width, conditional formatting and selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Unmerge the ranges that are going to change shape before editing values
# ---------------------------------------------------------------------------
$ws.Range("A2:A7").UnMerge()
$ws.Range("B2:B6").UnMerge()

# ---------------------------------------------------------------------------
# 2. Update existing rows 2-7
# ---------------------------------------------------------------------------

# Row 2 - Cantidad de registros migrados - Pago Unico (Terminado)
$ws.Range("D2").Value = "Creacion de los PR para la verificacion de cantidades de registros migrados de Pago Unico"
$ws.Range("G2").Value = 42223
$ws.Range("H2").Value = 42166
$ws.Range("I2").Value = 1
$ws.Range("K2").Value = "Terminado"
$ws.Range("M2").Value = "Script SQL desarrollado"

# Row 3 - Cantidad de registros migrados - Reprocesos (Terminado)
$ws.Range("D3").Value = "Creacion de los PR para la verificacion de cantidades de registros migrados de Reprocesos"
$ws.Range("G3").Value = 42223
$ws.Range("H3").Value = 42166
$ws.Range("I3").Value = 1
$ws.Range("K3").Value = "Terminado"
$ws.Range("M3").Value = "Script SQL desarrollado"

# Row 4 - Cantidad de registros migrados - Workflow (En proceso)
$ws.Range("D4").Value = "Creacion de los PR para la verificacion de cantidades de registros migrados de Workflow"
$ws.Range("G4").Value = 42223
$ws.Range("H4").Value = 42166
$ws.Range("I4").Value = 0.2
$ws.Range("K4").Value = "En proceso"
$ws.Range("M4").Value = "Script SQL desarrollado"

# Row 5 - Dato a dato - Pago Unico (Inicial)
$ws.Range("D5").Value = "Creacion de los PR para la verificacion de registros migrados dato a dato de Pago Unico"
$ws.Range("G5").Value = 42228
$ws.Range("H5").Value = 42237
$ws.Range("I5").Value = 0.05
$ws.Range("K5").Value = "Inicial"
$ws.Range("M5").Value = "Script SQL desarrollado"

# Row 6 - Dato a dato - Reprocesos (Inicial)
$ws.Range("D6").Value = "Creacion de los PR para la verificacion de registros migrados dato a dato de Reprocesos"
$ws.Range("F6").Value = "Media"
$ws.Range("G6").Value = 42228
$ws.Range("H6").Value = 42237
$ws.Range("I6").Value = 0.05
$ws.Range("K6").Value = "Inicial"
$ws.Range("M6").Value = "Script SQL desarrollado"

# Row 7 - Dato a dato - Workflow (Inicial)
$ws.Range("D7").Value = "Creacion de los PR para la verificacion de registros migrados dato a dato de Workflow"
$ws.Range("F7").Value = "Media"
$ws.Range("G7").Value = 42228
$ws.Range("H7").Value = 42237
$ws.Range("I7").Value = 0.05
$ws.Range("K7").Value = "Inicial"
$ws.Range("M7").Value = "Script SQL desarrollado"

# Fix up the style of M7 (no longer the last row so it loses the special
# bottom border) and the A7/B7 cells, which now sit in the middle / bottom
# of the (growing) merged blocks. Grab the "middle" (A3, s=13) and "bottom"
# (B6, s=14) template styles *before* touching A6/A7/B6/B7 below.
$ws.Range("A3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # stash "middle" style (s=13) in scratch cell Z1
$ws.Range("B6").Copy()
$ws.Range("Z2").PasteSpecial(-4122)   # stash "bottom" style (s=14) in scratch cell Z2
$ws.Range("M2").Copy()
$ws.Range("Z3").PasteSpecial(-4122)   # stash "normal data" M-style (s=7) in scratch cell Z3

$ws.Range("Z1").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("Z1").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("M7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Add the two new rows (8 and 9)
# ---------------------------------------------------------------------------

# Reuse the column C:M formatting already used by the other data rows
$ws.Range("C5:M5").Copy()
$ws.Range("C8:M8").PasteSpecial(-4122)
$ws.Range("C5:M5").Copy()
$ws.Range("C9:M9").PasteSpecial(-4122)

$ws.Range("Z1").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("A9").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("Z2").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("Z3").Copy()
$ws.Range("M8").PasteSpecial(-4122)
$ws.Range("Z3").Copy()
$ws.Range("M9").PasteSpecial(-4122)

# Clean up the scratch cells used as style templates
$ws.Range("Z1:Z3").Clear()

$ws.Range("C8").Value = 7
$ws.Range("D8").Value = "Correr los PR de control de calidad "
$ws.Range("E8").Value = "Equipo"
$ws.Range("F8").Value = "Alta"
$ws.Range("G8").Value = 42240
$ws.Range("H8").Value = 42243
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = "Inicial"
$ws.Range("B8").Value = "Verificacion del Control de Calidad"

$ws.Range("C9").Value = 8
$ws.Range("D9").Value = "Revision de resultados"
$ws.Range("E9").Value = "Equipo"
$ws.Range("F9").Value = "Alta"
$ws.Range("G9").Value = 42244
$ws.Range("H9").Value = 42247
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = "Inicial"

# ---------------------------------------------------------------------------
# 4. Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 49.5
$ws.Rows.Item(3).RowHeight = 49.5
$ws.Rows.Item(4).RowHeight = 36
$ws.Rows.Item(5).RowHeight = 36
$ws.Rows.Item(6).RowHeight = 37.5
$ws.Rows.Item(7).RowHeight = 39.75
$ws.Rows.Item(8).RowHeight = 20.25
$ws.Rows.Item(9).RowHeight = 21.75

# ---------------------------------------------------------------------------
# 5. Column F width
# ---------------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 12.14

# ---------------------------------------------------------------------------
# 6. Merge cells to their new shape
# ---------------------------------------------------------------------------
$ws.Range("A2:A9").Merge()
$ws.Range("B2:B7").Merge()
$ws.Range("B8:B9").Merge()

# ---------------------------------------------------------------------------
# 7. Conditional formatting - rebuild to match the new layout
# ---------------------------------------------------------------------------
$ws.Cells.FormatConditions.Delete()

function Add-NuevoRule($rng, $col) {
    $rule = $rng.FormatConditions.Add(9, 0, "Nuevo")
    $rule.Formula1 = '=NOT(ISERROR(SEARCH("Nuevo",' + $col + ')))'
    $rule.Text = "Nuevo"
    $rule.Interior.Color = 49407
}

Add-NuevoRule $ws.Range("H2") "H2"
Add-NuevoRule $ws.Range("G2") "G2"
Add-NuevoRule $ws.Range("H3:H4") "H3"
Add-NuevoRule $ws.Range("G3:G4") "G3"
Add-NuevoRule $ws.Range("H5") "H5"
Add-NuevoRule $ws.Range("G5") "G5"
Add-NuevoRule $ws.Range("H6:H7") "H6"
Add-NuevoRule $ws.Range("G6:G7") "G6"
Add-NuevoRule $ws.Range("H8") "H8"
Add-NuevoRule $ws.Range("G8") "G8"
Add-NuevoRule $ws.Range("H9") "H9"
Add-NuevoRule $ws.Range("G9") "G9"

# ---------------------------------------------------------------------------
# 8. Selection
# ---------------------------------------------------------------------------
$ws.Range("J7").Select()
